$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '247.09'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.26%'
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '30.12'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '10.78%'
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.190'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.54%'
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05740'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.12%'
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.591'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.67%'
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8589'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '5.20%'
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8725'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '4.95%'
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1363'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.45%'
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07147'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '3.44%'
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.25%'
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09397'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.10%'
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001522'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.79%'
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.04140'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-7.80%'
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006109'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.12%'
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.506'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.81%'
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.044'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.77%'
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.276'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.32%'
$ws.Range("E18").ClearFormats()

# Row 19
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.01020'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1,615.93%'
$ws.Range("E19").ClearFormats()

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.32%'
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.03288'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '5.98%'
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.24%'
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.603'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-3.75%'
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.44%'
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001213'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.97%'
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.47%'
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '20.45%'
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '-0.29%'
$ws.Range("E28").ClearFormats()

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005707'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-5.90%'
$ws.Range("E41").ClearFormats()

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002199'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-15.68%'
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009458'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '15.25%'
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005103'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-3.81%'
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.03%'
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.08895'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-18.33%'
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002766'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '6.57%'
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.03%'
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.03%'
$ws.Range("E50").ClearFormats()
